$d = $word.ActiveDocument

# --- Change 1: merge "Four fuzzing products..." sentence, dropping the
#     proofErr spell-check wrappers around Sulley/Codenomicon.
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute(
    "Four fuzzing products were included in the analysis: two open source solutions (Sulley and Peach Community) and two commercial solutions (Peach Professional and Codenomicon ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Four fuzzing products were included in the analysis: two open source solutions (Sulley and Peach Community) and two commercial solutions (Peach Professional and Codenomicon ",
    2)

# --- Change 2: merge "The team talked with sales staff..." sentence,
#     dropping the proofErr wrappers around Codenomicon/Fuzzer. (Leave the
#     following " Both products are licensed..." run untouched.)
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute(
    "The team talked with sales staff of both Codenomicon and Peach Fuzzer for budgetary, ballpark cost estimates for their products.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The team talked with sales staff of both Codenomicon and Peach Fuzzer for budgetary, ballpark cost estimates for their products.",
    2)

Write-Output "done"
